$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last refreshed" timestamp shown in A1.
$ws.Range("A1").Value = "Datos actualizados a 29 de Marzo de 2020 a las 21:20"

# Each entry: row, country, Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes
$rows = @(
    @(4, "Estados Unidos", 135957, 12379, 4378, 129188, 2948, 171, 2391),
    @(8, "Alemania", 61164, 3469, 9211, 51463, 1979, 57, 490),
    @(18, "Canada", 6258, 603, 508, 5687, 120, 3, 63),
    @(20, "Noruega", 4264, 249, 7, 4232, 91, 2, 25),
    @(34, "Rumania", 1760, 308, 206, 1511, 31, 6, 43),
    @(41, "Sudafrica", 1280, 93, 31, 1247, 7, 1, 2),
    @(42, "Finlandia", 1240, 73, 10, 1219, 32, 2, 11),
    @(77, "Tunez", 312, 34, 2, 302, 10, 0, 8),
    @(78, "Uruguay", 304, 0, 0, 303, 9, 0, 1),
    @(79, "Taiwan", 298, 15, 39, 257, 0, 0, 2),
    @(80, "Costa Rica", 295, 0, 3, 290, 6, 0, 2),
    @(146, "Republica de Yibuti", 18, 4, 0, 18, 0, 0, 0),
    @(148, "Mali", 18, 0, 0, 17, 0, 0, 1),
    @(149, "Islas Virgenes de los Estados Unidos", 17, 0, 0, 17, 0, 0, 0),
    @(150, "Maldivas", 17, 1, 13, 4, 0, 0, 0),
    @(151, "Guinea", 16, 8, 0, 16, 0, 0, 0),
    @(152, "Nueva Caledonia", 15, 0, 0, 15, 0, 0, 0),
    @(153, "Haiti", 15, 7, 1, 14, 0, 0, 0),
    @(157, "San Martin (Parte Francesa)", 11, 0, 0, 11, 0, 0, 0),
    @(158, "Dominica", 11, 0, 0, 11, 0, 0, 0)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]
    $ws.Cells.Item($rowNum, 7).Value = $r[7]
    $ws.Cells.Item($rowNum, 8).Value = $r[8]
}
